$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (kept as text to match original inlineStr formatting)
$updates = @{
    "D2" = "317.78"
    "E2" = "2.53%"
    "D3" = "41.36"
    "E3" = "1.32%"
    "D4" = "5.227"
    "E4" = "2.21%"
    "D5" = "0.07667"
    "E5" = "0.01%"
    "D6" = "1.672"
    "E6" = "3.66%"
    "D7" = "0.9349"
    "E7" = "2.23%"
    "D9" = "0.1270"
    "E9" = "0.75%"
    "D10" = "0.1840"
    "E10" = "1.86%"
    "D11" = "0.09151"
    "E11" = "0.95%"
    "D12" = "0.04128"
    "E12" = "-3.79%"
    "D13" = "0.1052"
    "E13" = "0.37%"
    "D14" = "0.001274"
    "E14" = "3.68%"
    "D15" = "0.005992"
    "E15" = "3.76%"
    "D17" = "3.345"
    "E17" = "-0.17%"
    "E18" = "1.75%"
    "D19" = "0.3348"
    "E19" = "0.93%"
    "D20" = "8.421"
    "E20" = "22.14%"
    "D21" = "0.1360"
    "E21" = "-2.18%"
    "D23" = "0.04043"
    "E23" = "0.36%"
    "D24" = "0.001276"
    "E24" = "0.74%"
    "D25" = "0.004083"
    "E25" = "-1.39%"
    "D26" = "0.0001280"
    "E26" = "0.90%"
    "D38" = "0.02510"
    "E38" = "3.06%"
    "D39" = "0.05264"
    "E39" = "0.05%"
    "D40" = "0.007804"
    "E40" = "-0.28%"
    "E41" = "-0.56%"
    "D42" = "0.007090"
    "E42" = "4.34%"
    "D43" = "0.002063"
    "E43" = "12.20%"
    "D44" = "0.008336"
    "E44" = "1.71%"
    "D45" = "0.3461"
    "E45" = "3.37%"
    "D46" = "0.00006673"
    "E46" = "-2.49%"
    "D47" = "0.00000000755"
    "E47" = "0.78%"
    "D48" = "0.1982"
    "E48" = "-37.93%"
    "D49" = "0.004227"
    "E49" = "40.79%"
    "D50" = "0.00002114"
    "E50" = "0.78%"
    "D51" = "0.0002013"
    "E51" = "0.78%"
}

# Force text number format first so numeric-looking strings are not coerced to numbers
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).NumberFormat = "@"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# Restore default (Normal) style so cells keep the original unstyled appearance
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Style = "Normal"
}
